$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '34.163.89'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.66%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.789.95'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.72%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '226.48'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.77%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.547'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.46%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '31.97'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.34%  '
$ws.Range("E9").Value = '  +1.57%  '
$ws.Range("E10").Value = '  -1.69%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0947'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.04%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '2.047.35'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.69%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '11.16'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +2.35%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.795.83'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.73%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '34.090.20'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.50%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.622'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.26%  '
$ws.Range("E17").Value = '  +1.15%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '68.14'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.71%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '245.19'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.11%  '
$ws.Range("E20").Value = '  -0.38%  '
$ws.Range("E22").Value = '  +1.48%  '
$ws.Range("E23").Value = '  +0.48%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.04'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.91%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '161.17'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.38%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.16'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.27%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '16.33'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.27%  '
$ws.Range("E28").Value = '  +0.84%  '
$ws.Range("E29").Value = '  +0.05%  '
$ws.Range("E30").Value = '  -0.03%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0519'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.77%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.67'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.62%  '
$ws.Range("E33").Value = '  +3.19%  '
$ws.Range("E34").Value = '  -0.07%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.455.01'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +4.55%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.649'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.88%  '
$ws.Range("E37").Value = '  +8.17%  '
$ws.Range("E38").Value = '  +2.66%  '
$ws.Range("E39").Value = '  -0.17%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '80.37'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +3.64%  '
$ws.Range("E41").Value = '  +0.46%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.921'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.32%  '
$ws.Range("E43").Value = '  +0.45%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '13.48'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.54%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0511'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.85%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '6.06'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +4.32%  '
$ws.Range("E47").Value = '  -0.01%  '
$ws.Range("E48").Value = '  -1.96%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.948.78'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.01%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '106.18'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.50%  '
$ws.Range("E51").Value = '  +0.02%  '
